$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.575.93"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").Value = "1.825.32"
$ws.Range("E3").Value = "  +0.86%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.007"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "308.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4683"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.98%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3592"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07135"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9280"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07663"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.76%  "

$ws.Range("D13").Value = "1.840.54"
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.251"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.337"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.009"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008534"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").Value = "26.607.88"
$ws.Range("E20").Value = "  +0.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.014"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.66%  "

$ws.Range("D23").Value = "2.075.84"
$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.911"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.26%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("E28").Value = "  -2.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.854"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08813"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.154"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.51%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.853"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.160"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7361"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.430"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.075"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01920"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.943"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05141"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.68%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.875"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5044"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1494"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.078"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.007"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4629"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.573"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06023"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.71%  "
